{"js": "// The author's commit (\"remove ep from has\") carves a new \"has_ep\"\n// tagged content control out of the plain descriptive sentence that used\n// to read:\n//   \"A gyomor, a nyomb\u00e9l \u00e9s a t\u00f6bbi b\u00e9lszakasz kp. t\u00e1g, ...\"\n// After the edit it reads:\n//   \"A \" + [has_ep content control containing a single space] + \"t\u00f6bbi b\u00e9lszakasz kp. t\u00e1g, ...\"\n// i.e. the words \"gyomor, a nyomb\u00e9l \u00e9s a \" are removed from the sentence\n// and replaced by a new (empty-looking) content control named \"has_ep\",\n// matching the style of the many other \"has_*\" content controls already\n// present in this document (e.g. has_lep, has_maj, has_vese).\n\nconst searchResults = context.document.body.search(\"gyomor, a nyomb\u00e9l \u00e9s a \", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Expected text 'gyomor, a nyomb\u00e9l \u00e9s a ' was not found in the document.\");\n}\n\nconst target = searchResults.items[0];\n\n// Wrap the words that are being dropped from the sentence in a brand new\n// content control first (this lets the control inherit the run formatting\n// of the wrapped text into its properties), tagged/aliased \"has_ep\" and\n// using the same \"tags\" appearance used by the sibling has_* controls.\nconst hasEpControl = target.insertContentControl();\nhasEpControl.tag = \"has_ep\";\nhasEpControl.title = \"has_ep\";\nhasEpControl.appearance = \"Tags\";\nawait context.sync();\n\n// Now collapse the content control's text down to the single space that\n// becomes its (empty-looking) content.\nhasEpControl.insertText(\" \", \"Replace\");\nawait context.sync();\n", "ps1": "# The author's commit (\"remove ep from has\") carves a new \"has_ep\"\n# tagged content control out of the plain descriptive sentence that used\n# to read:\n#   \"A gyomor, a nyomb\u00e9l \u00e9s a t\u00f6bbi b\u00e9lszakasz kp. t\u00e1g, ...\"\n# After the edit it reads:\n#   \"A \" + [has_ep content control containing a single space] + \"t\u00f6bbi b\u00e9lszakasz kp. t\u00e1g, ...\"\n# i.e. the words \"gyomor, a nyomb\u00e9l \u00e9s a \" are removed from the sentence\n# and replaced by a new (empty-looking) content control named \"has_ep\",\n# matching the style of the many other \"has_*\" content controls already\n# present in this document (e.g. has_lep, has_maj, has_vese).\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"gyomor, a nyomb\u00e9l \u00e9s a \")\nif (-not $found) {\n    throw \"Expected text 'gyomor, a nyomb\u00e9l \u00e9s a ' was not found in the document.\"\n}\n\n# Wrap the words that are being dropped from the sentence in a brand new\n# content control first (this lets the control inherit the run formatting\n# of the wrapped text into its properties), tagged/aliased \"has_ep\" and\n# using the same \"tags\" appearance used by the sibling has_* controls.\n# wdContentControlRichText = 0\n$cc = $d.ContentControls.Add(0, $rng)\n$cc.Tag = \"has_ep\"\n$cc.Title = \"has_ep\"\n$cc.Appearance = \"Tags\"\n\n# Now collapse the content control's text down to the single space that\n# becomes its (empty-looking) content.\n$cc.Range.Text = \" \"\n"}
